$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3511.8096
$ws.Range("J76").Value = 4352.8887
$ws.Range("L76").Value = 4352.8887
$ws.Range("N76").Value = -4982.8887
$ws.Range("H79").Value = 3511.8096
$ws.Range("J79").Value = 4352.8887
$ws.Range("L79").Value = 4352.8887
$ws.Range("N79").Value = -6536.8887
$ws.Range("H96").Value = 347.22223
$ws.Range("I96").Value = 219.2
$ws.Range("J96").Value = 507.25
$ws.Range("K96").Value = 657.5999999999999
$ws.Range("L96").Value = 1521.75
$ws.Range("M96").Value = 715.4000000000001
$ws.Range("N96").Value = -4267.75
$ws.Range("H100").Value = 33335192
$ws.Range("I100").Value = 41668140
$ws.Range("J100").Value = 3400
$ws.Range("K100").Value = 41668140
$ws.Range("L100").Value = 3400
$ws.Range("M100").Value = -41667599
$ws.Range("N100").Value = -4482
$ws.Range("H113").Value = 7805.0586
$ws.Range("I113").Value = 6656.2856
$ws.Range("J113").Value = 8609.200000000001
$ws.Range("K113").Value = 6656.2856
$ws.Range("L113").Value = 8609.200000000001
$ws.Range("M113").Value = -3402.2856
$ws.Range("N113").Value = -15117.2
$ws.Range("H133").Value = 25000
$ws.Range("J133").Value = 25000
$ws.Range("L133").Value = 25000
$ws.Range("N133").Value = -35120

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 10000
$ws.Range("J8").Value = 10000
$ws.Range("L8").Value = 10000
$ws.Range("N8").Value = -10288
$ws.Range("H32").Value = 4529.1445
$ws.Range("I32").Value = 4145.5
$ws.Range("J32").Value = 6303.5
$ws.Range("K32").Value = 4145.5
$ws.Range("L32").Value = 6303.5
$ws.Range("M32").Value = -3858.5
$ws.Range("N32").Value = -6877.5
$ws.Range("H110").Value = 679.13336
$ws.Range("I110").Value = 599
$ws.Range("J110").Value = 1200
$ws.Range("K110").Value = 599
$ws.Range("L110").Value = 1200
$ws.Range("M110").Value = 1446
$ws.Range("N110").Value = -5290
$ws.Range("H132").Value = 2152.5593
$ws.Range("I132").Value = 1243.921
$ws.Range("J132").Value = 3796.762
$ws.Range("K132").Value = 3731.763
$ws.Range("L132").Value = 11390.286
$ws.Range("M132").Value = -1201.763
$ws.Range("N132").Value = -16450.286
$ws.Range("H133").Value = 41464.2
$ws.Range("J133").Value = 41464.2
$ws.Range("L133").Value = 41464.2
$ws.Range("N133").Value = -46524.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2586588.2
$ws.Range("I7").Value = 3606180.5
$ws.Range("J7").Value = 37607
$ws.Range("K7").Value = 3606180.5
$ws.Range("L7").Value = 37607
$ws.Range("M7").Value = -3606067.5
$ws.Range("N7").Value = -37833
$ws.Range("H105").Value = 2336.8667
$ws.Range("I105").Value = 1991
$ws.Range("J105").Value = 3028.6
$ws.Range("K105").Value = 1991
$ws.Range("L105").Value = 3028.6
$ws.Range("M105").Value = -244
$ws.Range("N105").Value = -6522.6
$ws.Range("H132").Value = 32000
$ws.Range("J132").Value = 32000
$ws.Range("L132").Value = 32000
$ws.Range("N132").Value = -42120
$ws.Range("H134").Value = 4163.478
$ws.Range("I134").Value = 6014.3477
$ws.Range("J134").Value = 2312.6086
$ws.Range("K134").Value = 18043.0431
$ws.Range("L134").Value = 6937.825800000001
$ws.Range("M134").Value = -15508.0431
$ws.Range("N134").Value = -12007.8258
$ws.Range("H139").Value = 28354.5
$ws.Range("I139").Value = 26709
$ws.Range("J139").Value = 30000
$ws.Range("K139").Value = 26709
$ws.Range("L139").Value = 30000
$ws.Range("M139").Value = -21569
$ws.Range("N139").Value = -40280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 25641624
$ws.Range("I16").Value = 25641624
$ws.Range("K16").Value = 25641624
$ws.Range("M16").Value = -25641337
$ws.Range("H113").Value = 25641624
$ws.Range("I113").Value = 25641624
$ws.Range("K113").Value = 25641624
$ws.Range("M113").Value = -25639454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 130592.875
$ws.Range("I5").Value = 12878.375
$ws.Range("J5").Value = 154739.44
$ws.Range("K5").Value = 38635.125
$ws.Range("L5").Value = 464218.32
$ws.Range("M5").Value = -38523.125
$ws.Range("N5").Value = -464442.32
$ws.Range("H68").Value = 2788.054
$ws.Range("I68").Value = 4181.2
$ws.Range("J68").Value = 1838.1818
$ws.Range("K68").Value = 12543.6
$ws.Range("L68").Value = 5514.5454
$ws.Range("M68").Value = -11732.6
$ws.Range("N68").Value = -7136.5454
$ws.Range("H71").Value = 2788.054
$ws.Range("I71").Value = 4181.2
$ws.Range("J71").Value = 1838.1818
$ws.Range("K71").Value = 37630.8
$ws.Range("L71").Value = 16543.6362
$ws.Range("M71").Value = -33574.8
$ws.Range("N71").Value = -24655.6362
$ws.Range("H134").Value = 10939.434
$ws.Range("I134").Value = 13131.444
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 39394.33199999999
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -34324.33199999999
$ws.Range("N134").Value = -40140
$ws.Range("H135").Value = 130592.875
$ws.Range("I135").Value = 12878.375
$ws.Range("J135").Value = 154739.44
$ws.Range("K135").Value = 115905.375
$ws.Range("L135").Value = 1392654.96
$ws.Range("M135").Value = -113370.375
$ws.Range("N135").Value = -1397724.96

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3340
$ws.Range("I102").Value = 1142.4
$ws.Range("J102").Value = 5171.3335
$ws.Range("K102").Value = 1142.4
$ws.Range("L102").Value = 5171.3335
$ws.Range("M102").Value = 479.5999999999999
$ws.Range("N102").Value = -8415.333500000001
$ws.Range("H113").Value = 55556760
$ws.Range("I113").Value = 90910050
$ws.Range("J113").Value = 1586.1428
$ws.Range("K113").Value = 90910050
$ws.Range("L113").Value = 1586.1428
$ws.Range("M113").Value = -90907880
$ws.Range("N113").Value = -5926.1428
$ws.Range("H122").Value = 10805469
$ws.Range("I122").Value = 12964962
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 38894886
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -38892436
$ws.Range("N122").Value = -28900
$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -50120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 15874337
$ws.Range("J22").Value = 1541.6666
$ws.Range("L22").Value = 1541.6666
$ws.Range("N22").Value = -2131.6666
$ws.Range("H27").Value = 15874337
$ws.Range("J27").Value = 1541.6666
$ws.Range("L27").Value = 1541.6666
$ws.Range("N27").Value = -1755.6666
$ws.Range("H46").Value = 27778988
$ws.Range("I46").Value = 55556540
$ws.Range("K46").Value = 55556540
$ws.Range("M46").Value = -55556352
$ws.Range("H61").Value = 2385.4285
$ws.Range("I61").Value = 1819.1
$ws.Range("J61").Value = 3801.25
$ws.Range("K61").Value = 1819.1
$ws.Range("L61").Value = 3801.25
$ws.Range("M61").Value = -1617.1
$ws.Range("N61").Value = -4205.25
$ws.Range("H82").Value = 108400.55
$ws.Range("I82").Value = 34633.332
$ws.Range("J82").Value = 136063.25
$ws.Range("K82").Value = 34633.332
$ws.Range("L82").Value = 136063.25
$ws.Range("M82").Value = -34272.332
$ws.Range("N82").Value = -136785.25
$ws.Range("H85").Value = 108400.55
$ws.Range("I85").Value = 34633.332
$ws.Range("J85").Value = 136063.25
$ws.Range("K85").Value = 34633.332
$ws.Range("L85").Value = 136063.25
$ws.Range("M85").Value = -33385.332
$ws.Range("N85").Value = -138559.25
$ws.Range("H113").Value = 2385.4285
$ws.Range("I113").Value = 1819.1
$ws.Range("J113").Value = 3801.25
$ws.Range("K113").Value = 1819.1
$ws.Range("L113").Value = 3801.25
$ws.Range("M113").Value = 350.9000000000001
$ws.Range("N113").Value = -8141.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 100000584
$ws.Range("I107").Value = 111111670
$ws.Range("K107").Value = 333335010
$ws.Range("M107").Value = -333333090
$ws.Range("H113").Value = 992.3570999999999
$ws.Range("I113").Value = 737.56525
$ws.Range("J113").Value = 2164.4
$ws.Range("K113").Value = 2212.69575
$ws.Range("L113").Value = 6493.200000000001
$ws.Range("M113").Value = -42.69574999999986
$ws.Range("N113").Value = -10833.2
$ws.Range("H132").Value = 1612.5333
$ws.Range("I132").Value = 1123.359
$ws.Range("J132").Value = 2521
$ws.Range("K132").Value = 3370.077
$ws.Range("L132").Value = 7563
$ws.Range("M132").Value = -840.0769999999998
$ws.Range("N132").Value = -12623
